# Daily market-data refresh: update watchlist prices/pct_change/rsi (sorted by
# rsi ascending), recompute the stocks sheet (current_price/value/performance/
# current_rsi pulled from the watchlist), roll the portfolio STOCKS/TOTAL
# forward, and append a new row to the summary log.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) watchlist sheet - new snapshot, already ordered by ascending RSI
# ---------------------------------------------------------------------------
$wl = $wb.Worksheets.Item("watchlist")

$watchlist = @(
    @("CSCO", 42.31, 1.0509, 18.8356164383562),
    @("XOM",  41.32, -1.5253, 34.77366255144038),
    @("DIS",  128.12, 0.2739, 35.63829787234049),
    @("IBM",  123.15, -0.5572, 41.64588528678311),
    @("CVX",  84.81, -1.8289, 43.67201426024958),
    @("GS",   201.85, -0.5763, 44.43929564411491),
    @("WBA",  39.71, -1.3416, 44.46153846153847),
    @("BA",   169.58, 0.1831, 45.65649867374009),
    @("TRV",  112.53, -0.5919, 45.79569145239752),
    @("KO",   47.35, -0.0422, 47.63636363636373),
    @("MSFT", 214.58, 2.3271, 47.78132765353214),
    @("UNH",  313.33, -0.6563, 48.72422229989511),
    @("JPM",  97.37, -1.1974, 50.47899778924101),
    @("AXP",  96.72, -0.1755, 52.79329608938547),
    @("INTC", 49.17, 1.7381, 54.83870967741937),
    @("WMT",  130.57, -1.3896, 55.77557755775576),
    @("PFE",  38.72, 1.2023, 58.96414342629474),
    @("RTX",  60.99, -1.0545, 59.8984771573604),
    @("CAT",  137.48, -0.3912, 61.64280331574982),
    @("AAPL", 473.1, 2.219, 63.26234970302764),
    @("DOW",  43.69, -2.1281, 64.80144404332133),
    @("JNJ",  151.42, 0.6849, 68.51674641148321),
    @("V",    204.15, 1.5722, 69.00937081659968),
    @("HD",   280.68, -0.7707, 69.73821989528794),
    @("VZ",   58.96, -0.1524, 70.84870848708491),
    @("PG",   136.85, 0.7955, 73.73853211009174),
    @("MMM",  162.97, -0.6099, 76.34164777021921),
    @("NKE",  108.01, -0.3506, 85.59322033898314),
    @("MCD",  209.88, 0.1766, 86.09913793103453),
    @("MRK",  85.03, 0, 91.63179916317982)
)

for ($i = 0; $i -lt $watchlist.Length; $i++) {
    $r = $i + 2
    $row = $watchlist[$i]
    $wl.Cells.Item($r, 1).Value = $row[0]
    $wl.Cells.Item($r, 2).Value = $row[1]
    $wl.Cells.Item($r, 3).Value = $row[2]
    $wl.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------------
# 2) stocks sheet - pull fresh current_price/current_rsi for each held
#    ticker from the watchlist snapshot above, recompute value/performance
# ---------------------------------------------------------------------------
$st = $wb.Worksheets.Item("stocks")

$stocks = @(
    @("XOM",  45.98, 41.32,  12, 495.84,            -10.1348, 34.77366255144038),
    @("INTC", 50.59, 49.17,  18, 885.0600000000001, -2.8069,  54.83870967741937),
    @("JNJ",  146.84, 151.42, 4, 605.6799999999999,  3.119,   68.51674641148321),
    @("TRV",  113.1, 112.53,  4, 450.12,             -0.504,  45.79569145239752),
    @("PFE",  37.79, 38.72,  14, 542.0799999999999,   2.461,  58.96414342629474),
    @("CSCO", 42.5,  42.31,  14, 592.34,             -0.4471, 18.8356164383562)
)

for ($i = 0; $i -lt $stocks.Length; $i++) {
    $r = $i + 2
    $row = $stocks[$i]
    $st.Cells.Item($r, 1).Value = $row[0]
    $st.Cells.Item($r, 2).Value = $row[1]
    $st.Cells.Item($r, 3).Value = $row[2]
    $st.Cells.Item($r, 4).Value = $row[3]
    $st.Cells.Item($r, 5).Value = $row[4]
    $st.Cells.Item($r, 6).Value = $row[5]
    $st.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 3) portfolio sheet - STOCKS rolls up the stocks value column, TOTAL is
#    CASH + STOCKS
# ---------------------------------------------------------------------------
$pf = $wb.Worksheets.Item("portfolio")

$newStocksValue = 3571.12
$newTotalValue = 11002.75

$pf.Cells.Item(3, 2).Value = $newStocksValue
$pf.Cells.Item(4, 2).Value = $newTotalValue

# ---------------------------------------------------------------------------
# 4) summary sheet - append today's snapshot as a new log row
# ---------------------------------------------------------------------------
$sm = $wb.Worksheets.Item("summary")

$newRow = 68
$prevRow = $newRow - 1

# Column A keeps the same bold/bordered/centered look as every other date
# cell in the log, so copy formatting from the row above before writing in
# the new timestamp.
$sm.Cells.Item($prevRow, 1).Copy()
$sm.Cells.Item($newRow, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$sm.Cells.Item($newRow, 1).Value = "20/08/2020 17:45:39"
$sm.Cells.Item($newRow, 2).Value = 7431.629999999999
$sm.Cells.Item($newRow, 3).Value = $newStocksValue
$sm.Cells.Item($newRow, 4).Value = $newTotalValue
